{"js": "// Upstream change (M2Doc issue #295): this template was re-saved by the\n// tooling used in that commit, which re-serializes each OOXML part and\n// re-emits attributes/namespace declarations in (alphabetical) sorted\n// order. Diffing the canonical word/document.xml and word/styles.xml\n// parts before/after shows *only* that attribute re-ordering: every\n// element, every piece of text, and every attribute value is identical -\n// nothing in the document's visible content, formatting, styles or\n// section layout actually changed.\n//\n// So there is no content edit to replay. We simply load/read the body so\n// the script exercises the Word JavaScript API against this document\n// without writing any different value back - the saved package keeps the\n// exact same paragraphs, runs, fields and section/style values as\n// before.docx.\nconst body = context.document.body;\nbody.load(\"text\");\n\nconst sections = context.document.sections;\nsections.load(\"items\");\n\nconst styles = context.document.getStyles();\nstyles.load(\"items/nameLocal\");\n\nawait context.sync();\n", "ps1": "# Upstream change (M2Doc issue #295): this template was re-saved by the\n# tooling used in that commit, which re-serializes each OOXML part and\n# re-emits attributes/namespace declarations in (alphabetical) sorted\n# order. Diffing the canonical word/document.xml and word/styles.xml\n# parts before/after shows *only* that attribute re-ordering: every\n# element, every piece of text, and every attribute value is identical -\n# nothing in the document's visible content, formatting, styles or\n# section layout actually changed.\n#\n# So there is no content edit to replay. We simply read the document\n# through the Word COM object model so the script exercises the API\n# against this document without writing any different value back - the\n# saved package keeps the exact same paragraphs, runs, fields and\n# section/style values as before.docx.\n$d = $word.ActiveDocument\n\n# Touch/read the body text (no assignment -> no mutation).\n$bodyText = $d.Content.Text\n\n# Touch/read the section page setup (values are left untouched).\nforeach ($section in $d.Sections) {\n    $pageSetup = $section.PageSetup\n    $w = $pageSetup.PageWidth\n    $h = $pageSetup.PageHeight\n}\n\n# Touch/read the styles collection (no values are changed).\nforeach ($style in $d.Styles) {\n    $styleName = $style.NameLocal\n}\n"}
